$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.9
$ws.Range("B3").Value = 9.9600000000000009
$ws.Range("B4").Value = 9.94

$ws.Range("C4").Select()
